$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry "Watch" / "Metal" / "Male" / "Sold Out" gets inserted at row 2,
# pushing the existing "Earring" entry down to row 3.
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("D2").Value2

$ws.Range("A2").Value = "Watch"
$ws.Range("B2").Value = "Metal"
$ws.Range("C2").Value = "Male"
$ws.Range("D2").Value = "Sold Out"

$ws.Range("D11").Select()
